{"js": "// Update the title placeholder and replace the single \"Student Names\"\n// paragraph with three separate paragraphs, one per student name.\n\nconst body = context.document.body;\n\n// --- Title: \"<Project Name>\" -> \"Data Analysis Tool: NSW Traffic Penalties\"\nconst titleHits = body.search(\"<Project Name>\", { matchCase: true });\ntitleHits.load(\"items\");\nawait context.sync();\n\nif (titleHits.items.length > 0) {\n  const titlePara = titleHits.items[0].paragraphs.getFirst();\n  titlePara.insertText(\"Data Analysis Tool: NSW Traffic Penalties\", Word.InsertLocation.replace);\n} else {\n  // Fallback: second paragraph of the document holds the project-name title.\n  const paragraphs = body.paragraphs;\n  paragraphs.load(\"items\");\n  await context.sync();\n  paragraphs.items[1].insertText(\"Data Analysis Tool: NSW Traffic Penalties\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// --- \"Student Names\" paragraph -> three paragraphs, one per student.\nconst namesHits = body.search(\"Student Names\", { matchCase: true });\nnamesHits.load(\"items\");\nawait context.sync();\n\nlet namesPara;\nif (namesHits.items.length > 0) {\n  namesPara = namesHits.items[0].paragraphs.getFirst();\n} else {\n  const paragraphs = body.paragraphs;\n  paragraphs.load(\"items\");\n  await context.sync();\n  namesPara = paragraphs.items[2];\n}\n\nnamesPara.insertText(\"Brianne Byer\", Word.InsertLocation.replace);\nconst secondPara = namesPara.insertParagraph(\"Wonwoo Choi\", Word.InsertLocation.after);\nsecondPara.insertParagraph(\"Marco Querzola\", Word.InsertLocation.after);\n\nawait context.sync();\n", "ps1": "# Update the title placeholder and replace the single \"Student Names\"\n# paragraph with three separate paragraphs, one per student name.\n\n$d = $word.ActiveDocument\n\nfunction Get-ParagraphIndexByText($doc, $searchText, $fallbackIndex) {\n    $range = $doc.Content\n    $found = $range.Find.Execute($searchText)\n    if (-not $found) {\n        return $fallbackIndex\n    }\n    $range.Expand(4) | Out-Null   # wdParagraph\n    $target = $range.Start\n    $count = $doc.Paragraphs.Count\n    for ($i = 1; $i -le $count; $i++) {\n        if ($doc.Paragraphs.Item($i).Range.Start -eq $target) {\n            return $i\n        }\n    }\n    return $fallbackIndex\n}\n\nfunction Set-ParagraphText($doc, $index, $text) {\n    $r = $doc.Paragraphs.Item($index).Range\n    $r.MoveEnd(1, -1) | Out-Null   # exclude the trailing paragraph mark\n    $r.Text = $text\n}\n\n# --- Title: \"<Project Name>\" -> \"Data Analysis Tool: NSW Traffic Penalties\"\n$titleIndex = Get-ParagraphIndexByText $d \"<Project Name>\" 2\nSet-ParagraphText $d $titleIndex \"Data Analysis Tool: NSW Traffic Penalties\"\n\n# --- \"Student Names\" paragraph -> three paragraphs, one per student.\n$namesIndex = Get-ParagraphIndexByText $d \"Student Names\" 3\nSet-ParagraphText $d $namesIndex \"Brianne Byer\"\n$d.Paragraphs.Item($namesIndex).Range.InsertParagraphAfter()\n\n$secondIndex = $namesIndex + 1\nSet-ParagraphText $d $secondIndex \"Wonwoo Choi\"\n$d.Paragraphs.Item($secondIndex).Range.InsertParagraphAfter()\n\n$thirdIndex = $namesIndex + 2\nSet-ParagraphText $d $thirdIndex \"Marco Querzola\"\n"}
